# Auto-generated: apply cryptos.xlsx value updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.326.04"
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("D3").Value = "2.244.52"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'237.32"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  -4.77%  "
$ws.Range("D7").Value = "'69.67"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -6.80%  "
$ws.Range("D10").Value = "'0.0994"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'58.89"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "'36.68"
$ws.Range("E12").Value = "  +14.14%  "
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'6.75"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").Value = "2.579.13"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "'15.05"
$ws.Range("E16").Value = "  -6.13%  "
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("D18").Value = "2.253.37"
$ws.Range("D19").Value = "42.239.74"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "'6.27"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").Value = "'73.35"
$ws.Range("D23").Value = "'234.49"
$ws.Range("E23").Value = "  -6.53%  "
$ws.Range("D24").Value = "'2.00"
$ws.Range("E24").Value = "  +5.59%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'3.67"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").Value = "'171.12"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").Value = "'0.126"
$ws.Range("E33").Value = "  -5.12%  "
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'5.32"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -6.85%  "
$ws.Range("D37").Value = "'3.73"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "'22.67"
$ws.Range("E38").Value = "  +21.63%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").Value = "'0.0276"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").Value = "'5.94"
$ws.Range("E41").Value = "  -6.51%  "
$ws.Range("D42").Value = "'65.32"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "'9.33"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "'4.94"
$ws.Range("E44").Value = "  -16.61%  "
$ws.Range("D45").Value = "'0.104"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  +13.50%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").Value = "'10.23"
$ws.Range("E49").Value = "  +10.74%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.19"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = "  -3.02%  "

Write-Host "Applied 85 cell updates"
